# Add a new "CUATRIMESTRE" column (I) with header + per-row semester values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header cell + formatting copied from the adjacent header (H1) so it matches
# the bold/filled/centered header style used across the table.
$ws.Cells.Item(1, 9).Value = "CUATRIMESTRE"
$ws.Cells.Item(1, 8).Copy() | Out-Null
$ws.Cells.Item(1, 9).PasteSpecial(-4122) | Out-Null

# Per-row CUATRIMESTRE values for data rows 2..71
$values = @(1,1,1,1,1,2,2,2,2,2,1,1,1,1,1,2,2,2,2,2,1,1,1,1,1,2,2,2,2,2,1,1,1,1,2,2,2,2,1,1,1,1,2,2,2,2,1,1,1,1,2,2,2,2,1,1,1,1,2,2,1,2,1,1,2,2,2,2,2,2)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i]
}

# Size the new column to fit its contents, like the other bestFit columns.
$ws.Columns.Item(9).ColumnWidth = 13.736979166666666

# Restore the selection/scroll state left in the sheet after editing.
$ws.Range("L58").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 32
$excel.ActiveWindow.ScrollColumn = 1
